# Insert a new weekly price-report row for "Coliflor" (Feria Lagunitas de
# Puerto Montt) right before the current row 353. This pushes the existing
# rows 353..381 down to 354..382 (and the used range grows to A1:R382),
# matching how a new day's entry gets prepended to this daily/weekly log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 353; everything below shifts down one row.
$ws.Rows.Item(353).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(353, 1).Value  = 4
$ws.Cells.Item(353, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(353, 3).Value  = "Los Lagos"
$ws.Cells.Item(353, 4).Value  = 44783
$ws.Cells.Item(353, 5).Value  = 10
$ws.Cells.Item(353, 6).Value  = 100112008
$ws.Cells.Item(353, 7).Value  = "Coliflor"
$ws.Cells.Item(353, 8).Value  = "Sin especificar"
$ws.Cells.Item(353, 9).Value  = "Primera"
$ws.Cells.Item(353, 10).Value = 100
$ws.Cells.Item(353, 11).Value = 1500
$ws.Cells.Item(353, 12).Value = 1500
$ws.Cells.Item(353, 13).Value = 1500
$ws.Cells.Item(353, 14).Value = "`$/unidad"
$ws.Cells.Item(353, 15).Value = "Región Metropolitana"
$ws.Cells.Item(353, 16).Value = 1500
$ws.Cells.Item(353, 17).Value = 1
$ws.Cells.Item(353, 18).Value = "Hortaliza"
